# Applies the "Updates xml for package update" edit to Location_attributes.xlsx
# - Adds two new attribute rows (Habitat Type, Purpose) to the "attribute" sheet
# - Extends the three list data-validations to cover the new rows
# - Moves the active-cell selection to E8

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

# --- Insert two new rows (5 and 6) by copying the formatting of row 3, which
#     already carries the style pattern we need: A -> s=4, B -> no style,
#     C/D/F -> s=1 (list-validated columns).
$ws1.Rows("3:3").Copy()
$ws1.Rows("5:5").Insert(-4121)   # xlShiftDown

$ws1.Rows("3:3").Copy()
$ws1.Rows("6:6").Insert(-4121)   # xlShiftDown

# Row 5/6 col B shouldn't carry any style (matches row 3's B column) - the
# insert operation above guesses a blended style for B, so reset it.
$ws1.Range("B5").ClearFormats()
$ws1.Range("B6").ClearFormats()

# --- Fill in the new attribute rows' values
$ws1.Range("A5").Value = "Habitat Type "
$ws1.Range("B5").Value = "Habitat type at each site"

$ws1.Range("A6").Value = "Purpose"
$ws1.Range("B6").Value = "Identified purpose of each site"

# --- Extend the data validation ranges from row 4 to row 6 (delete + re-add
#     in the original C, D, F order so the sqref ranges grow to match).
$ws1.Range("C2:C6").Validation.Delete()
$ws1.Range("D2:D6").Validation.Delete()
$ws1.Range("F2:F6").Validation.Delete()

$ws1.Range("C2:C6").Validation.Add(3, 1, 1, '"string,boolean,decimal,float,double,duration,dateTime,time,date,gYearMonth,gYear,gMonthDay,gDay,gMonth"')
$ws1.Range("D2:D6").Validation.Add(3, 1, 1, '"nominal,ordinal,interval,ratio,dateTime"')
$ws1.Range("F2:F6").Validation.Add(3, 1, 1, '"text,enumerated,dateTime,numeric"')

# --- Update the active selection shown when the sheet is next opened
$ws1.Range("E8").Select() | Out-Null
